# Add "Save" column (H) to the s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting from the neighboring header cell (G1),
# which carries the bold/centered/bordered header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for H2:H9 (plain/default formatting, matching columns B:G)
$values = @(1, 0, 0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
